# Build script at 2023-04-12 14:53:07 UTC
# Fills in the previously-empty Portuguese/English content cells for LOT2008
# (Objetivos, Docentes responsaveis, Programa resumido, Programa, Metodo,
# Criterio, Norma de recuperacao, Bibliografia), which shifts every row
# from 13 onward down by one and appends a new trailing row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; this pushes the old rows 13-23 down to
# 14-24 (carrying their formatting/row-heights with them) and opens up a
# blank row 13 for the "Docentes responsaveis" value.
$ws.Rows.Item(13).Insert()

# The insert leaves a bare formatted-but-empty A13 cell behind; drop it so
# the row only contains the B13/C13 values, as in the final layout.
$ws.Range("A13").Clear()

$ws.Range("B13").Value = "5840494 - Maria Eleonora Andrade de Carvalho"
$ws.Range("C13").Value = "5840494 - Maria Eleonora Andrade de Carvalho"
$ws.Range("B14").Value = "Introdução ao metabolismo, BioenergéticaOxidações biológicas, TransporteGlicídios - metabolismoFotossínteseLipídeos - metabolismoAminoácidos - metabolismoIntegração MetabólicaCiclos vitais: oxigênio, carbono, nitrogênio e enxofre"
$ws.Range("C14").Value = "Introdução ao metabolismo, BioenergéticaOxidações biológicas, TransporteGlicídios - metabolismoFotossínteseLipídeos - metabolismoAminoácidos - metabolismoIntegração MetabólicaCiclos vitais: oxigênio, carbono, nitrogênio e enxofre"
$ws.Range("B15").Value = "Introduction to metabolism. Bioenergetics. Biological oxidations. Glycides transport  metabolism,  photosynthesis,  lipid metabolism, metabolism of nitrogenous compounds: amino acids, integration and control of metabolic processes, vital cycles: oxygen, carbon, nitrogen and sulfur."
$ws.Range("C15").Value = "Introduction to metabolism. Bioenergetics. Biological oxidations. Glycides transport  metabolism,  photosynthesis,  lipid metabolism, metabolism of nitrogenous compounds: amino acids, integration and control of metabolic processes, vital cycles: oxygen, carbon, nitrogen and sulfur."
$ws.Range("B16").Value = "Introdução ao metabolismo. Proteínas, polissacarídeos, lipídios: vias catabólicas e anabólicas. Bioenergética. Variação de energia livre: relação com a constante de equilíbrio e com o potencial redox. Processos exergônicos. Papel do fosfato: potencial de transferência de grupo fosfato. Importância energética do ATP.Oxidações biológicas. Coenzimas transportadoras de prótons e elétrons: nucleotídeos, flavino nucleotídeos, coenzima Q. Desidrogenase piridino e flavino nucleotídeos dependentes. Oxidases. Estrutura da membrana mitocondrial. Cadeia respiratória: função. Fosforilação oxidativa.Transporte. Composição das membranas biológicas: constituição química, caráter ?barreira permeabilidade?. Carreadores e canais, ionóforos. Transporte: mediado e não mediado. Glicídios - metabolismo. Degradação anaeróbica e aeróbica de glicídios: glicólise - localização das enzimas operantes, reações, produção de NADH.H+, fosforilação ao nível de substrato, balanço energético; ciclo de Krebs - localização das enzimas operantes, reações, produção de coenzimas reduzidas, balanço energético. ?Shunt? das hexoses-fosfato (ciclo das pentoses): localização das enzimas operantes, reações (fases oxidativa e não oxidativa), produção de NADPH.H  (implicação fisiológica). Fermentações: definição, fermentação e respiração, matérias primas usadas em fermentação amilácea e sacarínea, agente de fermentação, fermentações anaeróbicas - alcoólica e lática, fermentações aeróbicas - acética e cítrica. Fotossíntese. Estrutura dos cloroplastos. Luz: energia eletromagnética. Papel da clorofila na fotossíntese. Fotofosforilação cíclica e não cíclica. Redução do NADP . Fotólise da água. Síntese do aceptor de CO2, Ru-1,5diP. Ciclo de Calvin.Lipídios - metabolismo. b-oxidação de ácidos graxos de cadeia: com número par de átomos de C, com número impar de átomos de C, ramificada, balanço energético da b-oxidação. a-oxidação e w-oxidação. Metabolismo do glicerol. Formação de corpos cetônicos. Biossíntese de ácidos graxos.Aminoácidos - vias catabólicas. Digestão de proteínas, enzimas envolvidas e zimogênios, absorção. Transaminação, desaminação oxidativa, aminoácidos cetogênicos e glicogênicos, descarboxilação, ciclo da uréia. Eliminação de nitrogênio, vertebrados ureotélicos.Integração metabólica. Interelação do metabolismo intermediário de glicídios, lipídios, aminoácidos e ácidos nucléicos. Metabólitos comuns ao metabolismo de glicídios, lipídios e aminoácidos.Ciclos vitais: oxigênio, carbono, nitrogênio e enxofre. Ciclos do oxigênio e do carbono. Ciclo do nitrogênio: fixação biológica, nitrificação, utilização do nitrato, incorporação de amônia em compostos orgânicos. Ciclo do enxofre: assimilação do sulfato"
$ws.Range("C16").Value = "Introdução ao metabolismo. Proteínas, polissacarídeos, lipídios: vias catabólicas e anabólicas. Bioenergética. Variação de energia livre: relação com a constante de equilíbrio e com o potencial redox. Processos exergônicos. Papel do fosfato: potencial de transferência de grupo fosfato. Importância energética do ATP.Oxidações biológicas. Coenzimas transportadoras de prótons e elétrons: nucleotídeos, flavino nucleotídeos, coenzima Q. Desidrogenase piridino e flavino nucleotídeos dependentes. Oxidases. Estrutura da membrana mitocondrial. Cadeia respiratória: função. Fosforilação oxidativa.Transporte. Composição das membranas biológicas: constituição química, caráter ?barreira permeabilidade?. Carreadores e canais, ionóforos. Transporte: mediado e não mediado. Glicídios - metabolismo. Degradação anaeróbica e aeróbica de glicídios: glicólise - localização das enzimas operantes, reações, produção de NADH.H+, fosforilação ao nível de substrato, balanço energético; ciclo de Krebs - localização das enzimas operantes, reações, produção de coenzimas reduzidas, balanço energético. ?Shunt? das hexoses-fosfato (ciclo das pentoses): localização das enzimas operantes, reações (fases oxidativa e não oxidativa), produção de NADPH.H  (implicação fisiológica). Fermentações: definição, fermentação e respiração, matérias primas usadas em fermentação amilácea e sacarínea, agente de fermentação, fermentações anaeróbicas - alcoólica e lática, fermentações aeróbicas - acética e cítrica. Fotossíntese. Estrutura dos cloroplastos. Luz: energia eletromagnética. Papel da clorofila na fotossíntese. Fotofosforilação cíclica e não cíclica. Redução do NADP . Fotólise da água. Síntese do aceptor de CO2, Ru-1,5diP. Ciclo de Calvin.Lipídios - metabolismo. b-oxidação de ácidos graxos de cadeia: com número par de átomos de C, com número impar de átomos de C, ramificada, balanço energético da b-oxidação. a-oxidação e w-oxidação. Metabolismo do glicerol. Formação de corpos cetônicos. Biossíntese de ácidos graxos.Aminoácidos - vias catabólicas. Digestão de proteínas, enzimas envolvidas e zimogênios, absorção. Transaminação, desaminação oxidativa, aminoácidos cetogênicos e glicogênicos, descarboxilação, ciclo da uréia. Eliminação de nitrogênio, vertebrados ureotélicos.Integração metabólica. Interelação do metabolismo intermediário de glicídios, lipídios, aminoácidos e ácidos nucléicos. Metabólitos comuns ao metabolismo de glicídios, lipídios e aminoácidos.Ciclos vitais: oxigênio, carbono, nitrogênio e enxofre. Ciclos do oxigênio e do carbono. Ciclo do nitrogênio: fixação biológica, nitrificação, utilização do nitrato, incorporação de amônia em compostos orgânicos. Ciclo do enxofre: assimilação do sulfato"
$ws.Range("B17").Value = "Introduction to metabolism. Proteins, polysaccharides, lipids.: catabolic and anabolic pathways. Bioenergetics. The direction of processes: free energy: reaction with the balance constant and with redox potential. Exergonic processes. Role of phosphate: potential of transference of phosphate group. Role of ATP as a free energy currency. Biological oxidations. Electron and proton transporters coenzymes: nucleotides, flavin nucleotides, coenzyme Q. Pyridine and flavin nucleotides- dehydrogenase dependent. Oxidases.Structure of mitochondrial membrane. Respiratory chain: function. Oxidative phosphorylation. Transport. Composition of biological membranes: chemical constitution, characterization, barrier, permeability. Carriers and canal ionophores. Transportation: mediated and non-mediated. Glycides  metabolism. Anaerobic and aerobic degradation of glycides: glycolysis  localization of enzymes, reactions, NaDH.H+ production, the first substrt level phosphorylation, energetic balance; Citric acid cycle  localization of operating enzymes, reactions, production of reduced coenzymes, energetic balance. the pentose phosphate pathway: localization of the enzymes, reactions (oxidative and non-oxidative phases), NaDH.H+ production (physiological implication). Fermentation: definition, fermentation and respiration, raw-materials used in the starch and sugar fermentations, metabolic fates of pyruvate: ethanol and lactic metabolism, acetic and citric. Photosynthesis. Chloroplasts structure. The light reactions. Role of the chlorophyll in the photosynthesis. Cyclic and non-cyclic phosphorylation. NADP reduction. Water photolysis. Synthesis of the acceptor of CO2, Ru-1, 5-diP. Calvin cycle. Lipids  metabolism. Beta-oxidation pathways, oxidation of fatty acids with odd-numbered carbon chains, energetic balance of beta-oxidation, alpha-oxidation and w-oxidation.  Glycerol metabolism. Formation of ketone bodies. Biosynthesis of fatty acids. Amino acids  catabolic pathways. Digestion of proteins, aspects of amino acid synthesis and degradation.Transamination, urea cycle. Metabolic integration. Common metabolites to the metabolism of glycides, lipids and amino acids. Vital cycles: oxygen, carbon, nitrogen and sulfur. Carbon and oxygen cycles. Nitrogen cycle: biological fixation, nitrification, use of nitrate, incorporation of ammonia in organic compounds. Sulfur cycle: sulfate assimilation."
$ws.Range("C17").Value = "Introduction to metabolism. Proteins, polysaccharides, lipids.: catabolic and anabolic pathways. Bioenergetics. The direction of processes: free energy: reaction with the balance constant and with redox potential. Exergonic processes. Role of phosphate: potential of transference of phosphate group. Role of ATP as a free energy currency. Biological oxidations. Electron and proton transporters coenzymes: nucleotides, flavin nucleotides, coenzyme Q. Pyridine and flavin nucleotides- dehydrogenase dependent. Oxidases.Structure of mitochondrial membrane. Respiratory chain: function. Oxidative phosphorylation. Transport. Composition of biological membranes: chemical constitution, characterization, barrier, permeability. Carriers and canal ionophores. Transportation: mediated and non-mediated. Glycides  metabolism. Anaerobic and aerobic degradation of glycides: glycolysis  localization of enzymes, reactions, NaDH.H+ production, the first substrt level phosphorylation, energetic balance; Citric acid cycle  localization of operating enzymes, reactions, production of reduced coenzymes, energetic balance. the pentose phosphate pathway: localization of the enzymes, reactions (oxidative and non-oxidative phases), NaDH.H+ production (physiological implication). Fermentation: definition, fermentation and respiration, raw-materials used in the starch and sugar fermentations, metabolic fates of pyruvate: ethanol and lactic metabolism, acetic and citric. Photosynthesis. Chloroplasts structure. The light reactions. Role of the chlorophyll in the photosynthesis. Cyclic and non-cyclic phosphorylation. NADP reduction. Water photolysis. Synthesis of the acceptor of CO2, Ru-1, 5-diP. Calvin cycle. Lipids  metabolism. Beta-oxidation pathways, oxidation of fatty acids with odd-numbered carbon chains, energetic balance of beta-oxidation, alpha-oxidation and w-oxidation.  Glycerol metabolism. Formation of ketone bodies. Biosynthesis of fatty acids. Amino acids  catabolic pathways. Digestion of proteins, aspects of amino acid synthesis and degradation.Transamination, urea cycle. Metabolic integration. Common metabolites to the metabolism of glycides, lipids and amino acids. Vital cycles: oxygen, carbon, nitrogen and sulfur. Carbon and oxygen cycles. Nitrogen cycle: biological fixation, nitrification, use of nitrate, incorporation of ammonia in organic compounds. Sulfur cycle: sulfate assimilation."
$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2"
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2"
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("B22").Value = "1.Nelson, D.L.; Cox, M.M. Lehninger Principles of Biochemistry. Third Edition, Worth Publisher, New York, 20002.Voet, D; Voet, J G.; Pratt, C.W. Fundamentos de Bioquímica. Editora ARTMED, Porto Alegre, 20003.Stryer, L. Biochemistry. W.H. Freeman Company, New York, 19884.Jain, M.K. Introduction to Biological Membranes. John Wiley & Sons Inc., New York, 1988"
$ws.Range("C22").Value = "1.Nelson, D.L.; Cox, M.M. Lehninger Principles of Biochemistry. Third Edition, Worth Publisher, New York, 20002.Voet, D; Voet, J G.; Pratt, C.W. Fundamentos de Bioquímica. Editora ARTMED, Porto Alegre, 20003.Stryer, L. Biochemistry. W.H. Freeman Company, New York, 19884.Jain, M.K. Introduction to Biological Membranes. John Wiley & Sons Inc., New York, 1988"
